# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.579.82"
$ws.Range("E2").Value = "  +5.89%  "
$ws.Range("D3").Value = "2.297.70"
$ws.Range("E3").Value = "  +3.27%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'304.50"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").Value = "'100.54"
$ws.Range("E6").Value = "  +11.36%  "
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  +6.51%  "
$ws.Range("D10").Value = "'36.47"
$ws.Range("E10").Value = "  +10.56%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").Value = "'7.41"
$ws.Range("E12").Value = "  +6.65%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "2.648.78"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").Value = "2.302.46"
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("D16").Value = "'13.84"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("D18").Value = "46.570.95"
$ws.Range("E18").Value = "  +6.21%  "
$ws.Range("D19").Value = "'13.05"
$ws.Range("E19").Value = "  +11.82%  "
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("D21").Value = "'6.02"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").Value = "'66.21"
$ws.Range("E22").Value = "  +3.16%  "
$ws.Range("D23").Value = "'248.79"
$ws.Range("E23").Value = "  +5.71%  "
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +3.59%  "
$ws.Range("D27").Value = "'42.64"
$ws.Range("E27").Value = "  +9.12%  "
$ws.Range("E28").Value = "  +2.68%  "
$ws.Range("D29").Value = "'9.89"
$ws.Range("E29").Value = "  +5.72%  "
$ws.Range("D30").Value = "'20.02"
$ws.Range("E30").Value = "  +4.47%  "
$ws.Range("E31").Value = "  +12.63%  "
$ws.Range("E32").Value = "  +4.63%  "
$ws.Range("D33").Value = "'147.94"
$ws.Range("E33").Value = "  -2.12%  "
$ws.Range("D34").Value = "'0.0795"
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("E35").Value = "  +16.05%  "
$ws.Range("E36").Value = "  +10.38%  "
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("E38").Value = "  +6.12%  "
$ws.Range("D39").Value = "'16.03"
$ws.Range("E39").Value = "  +20.45%  "
$ws.Range("D40").Value = "'3.99"
$ws.Range("E40").Value = "  +10.99%  "
$ws.Range("D41").Value = "'3.35"
$ws.Range("E41").Value = "  +6.04%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'1.98"
$ws.Range("E44").Value = "  +11.01%  "
$ws.Range("D45").Value = "1.820.31"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").Value = "'88.21"
$ws.Range("E46").Value = "  +21.30%  "
$ws.Range("E47").Value = "  +6.37%  "
$ws.Range("D48").Value = "'73.54"
$ws.Range("E48").Value = "  +8.48%  "
$ws.Range("E49").Value = "  +6.66%  "
$ws.Range("D50").Value = "'96.07"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.525.27"
$ws.Range("E51").Value = "  +3.32%  "
